$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Día" / "Cantidad" data set (rows 2-12), replacing the old rows 2-10.
# The invalid "Invalid Date" row is gone; two extra dated rows were added
# along with a running count in column B.
$days = @(
  "Tue Nov 15 2022",
  "Mon Oct 31 2022",
  "Sat Oct 29 2022",
  "Thu Oct 27 2022",
  "Tue Nov 22 2022",
  "Tue Oct 25 2022",
  "Wed Oct 26 2022",
  "Thu Nov 17 2022",
  "Wed Nov 02 2022",
  "Fri Nov 18 2022",
  "Mon Oct 24 2022"
)

for ($i = 0; $i -lt $days.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $days[$i]
  $ws.Cells.Item($row, 2).Value = $i
}
